$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly report (2 rows: "Primera" / "Segunda" quality grades for
# Betarraga) is inserted at the top of the data block, right before the
# existing row 150. Everything that used to occupy rows 150-269 shifts
# down two rows, to 152-271, growing the used range from A1:R269 to
# A1:R271.
$ws.Rows("150:151").Insert()

# Seed the two new rows by duplicating the row that lands right below them
# (which, post-insert, holds the data that used to be row 150/151) so all
# the shared/common fields (market, region, category, unit, origin, etc.)
# come along for free.
$ws.Range("A152:R152").Copy()
$ws.Range("A150:R150").PasteSpecial()

$ws.Range("A153:R153").Copy()
$ws.Range("A151:R151").PasteSpecial()

# New row 150 - "Primera" quality for the new reporting date.
$ws.Cells.Item(150, 4).Value2 = 44658
$ws.Cells.Item(150, 11).Value2 = 450
$ws.Cells.Item(150, 12).Value2 = 500
$ws.Cells.Item(150, 13).Value2 = 475
$ws.Cells.Item(150, 16).Value2 = 158

# New row 151 - "Segunda" quality for the new reporting date.
$ws.Cells.Item(151, 4).Value2 = 44658
$ws.Cells.Item(151, 11).Value2 = 350
$ws.Cells.Item(151, 12).Value2 = 400
$ws.Cells.Item(151, 13).Value2 = 375
$ws.Cells.Item(151, 16).Value2 = 125
